$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "small tylostyle "
$ws.Range("A2").Value = "large tylostyle "
